$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto Price (column D) and Volume(1h) (column E) values.
# Some Price values are plain decimals ("0.481", "17.87", ...).
# Excel.Range.Value would auto-convert those into numbers, but the
# source column stores them as plain text, so we briefly force a
# text number format, assign the text, then restore the original
# cell style so formatting is left exactly as it was.

$ws.Range('D2').Value = '25.937.30'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.591.87'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.14%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.19'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('E6').Value = '  +0.09%  '
$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.481'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('E9').Value = '  -1.44%  '
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.87'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  -2.08%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0808'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  +2.53%  '
$ws.Range('D12').Value = '1.815.26'
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('D13').Value = '1.601.84'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('E14').Value = '  -1.07%  '
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.510'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').Value = '25.947.86'
$ws.Range('E16').Value = '  +0.20%  '
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '59.92'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('E19').Value = '  +0.19%  '
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '198.93'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +2.53%  '
$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.21'
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  +0.19%  '
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.21'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  -2.15%  '
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.98'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  +0.66%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.79'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +5.26%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.90'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -8.47%  '
$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.03'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  -0.84%  '
$ws.Range('E29').Value = '  -0.67%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('E32').Value = '  -0.58%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.93'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -3.61%  '
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('E35').Value = '  +1.79%  '
$ws.Range('D36').Value = '1.124.77'
$ws.Range('E36').Value = '  +1.62%  '
$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0162'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  +7.97%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('E39').Value = '  -0.96%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.781'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('E41').Value = '  -3.66%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.781'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  -2.85%  '
$ws.Range('D43').Value = '1.726.48'
$ws.Range('E43').Value = '  +0.23%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '92.38'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  -1.01%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.07'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('E46').Value = '  -1.78%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '53.20'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('E48').Value = '  -1.68%  '
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.407'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('E50').Value = '  +0.31%  '
$ws.Range('E51').Value = '  -17.53%  '
